$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws2 = $wb.Worksheets.Item("ARM")
$ws3 = $wb.Worksheets.Item("BSM")
$ws4 = $wb.Worksheets.Item("CRP")
$ws5 = $wb.Worksheets.Item("CUL")
$ws6 = $wb.Worksheets.Item("GSM")
$ws7 = $wb.Worksheets.Item("LTW")
$ws8 = $wb.Worksheets.Item("WVR")

# --- ALC ---
$ws1.Range("H112").Value = 2198.125
$ws1.Range("J112").Value = 3397.6
$ws1.Range("L112").Value = 10192.8
$ws1.Range("N112").Value = -12408.8
$ws1.Range("H113").Value = 1685
$ws1.Range("I113").Value = 1404.6666
$ws1.Range("J113").Value = 1778.4445
$ws1.Range("K113").Value = 1404.6666
$ws1.Range("L113").Value = 1778.4445
$ws1.Range("M113").Value = 1849.3334
$ws1.Range("N113").Value = -8286.4445
$ws1.Range("H115").Value = 100000000
$ws1.Range("I115").Value = 100000000
$ws1.Range("K115").Value = 300000000
$ws1.Range("M115").Value = -299998433
$ws1.Range("H116").Value = 2201058.8
$ws1.Range("I116").Value = 7695047
$ws1.Range("J116").Value = 3463.28
$ws1.Range("K116").Value = 7695047
$ws1.Range("L116").Value = 3463.28
$ws1.Range("M116").Value = -7691605
$ws1.Range("N116").Value = -10347.28
$ws1.Range("H118").Value = 559.75
$ws1.Range("I118").Value = 559.75
$ws1.Range("J118").Value = 0
$ws1.Range("K118").Value = 1679.25
$ws1.Range("L118").Value = 0
$ws1.Range("M118").Value = -22.25
$ws1.Range("N118").ClearContents()
$ws1.Range("H132").Value = 2696.138
$ws1.Range("I132").Value = 2799.4783
$ws1.Range("J132").Value = 2300
$ws1.Range("K132").Value = 8398.4349
$ws1.Range("L132").Value = 6900
$ws1.Range("M132").Value = -5868.4349
$ws1.Range("N132").Value = -11960
$ws1.Range("H136").Value = 50019.023
$ws1.Range("I136").Value = 50000
$ws1.Range("J136").Value = 50019.5
$ws1.Range("K136").Value = 50000
$ws1.Range("L136").Value = 50019.5
$ws1.Range("M136").Value = -44900
$ws1.Range("N136").Value = -60219.5

# --- ARM ---
$ws2.Range("H4").Value = 299
$ws2.Range("I4").Value = 298.5
$ws2.Range("J4").Value = 300.5
$ws2.Range("K4").Value = 298.5
$ws2.Range("L4").Value = 300.5
$ws2.Range("M4").Value = -182.5
$ws2.Range("N4").Value = -532.5
$ws2.Range("H5").Value = 262
$ws2.Range("I5").Value = 190
$ws2.Range("J5").Value = 550
$ws2.Range("K5").Value = 190
$ws2.Range("L5").Value = 550
$ws2.Range("M5").Value = -78
$ws2.Range("N5").Value = -774
$ws2.Range("H113").Value = 67448.5
$ws2.Range("J113").Value = 67448.5
$ws2.Range("L113").Value = 67448.5
$ws2.Range("N113").Value = -76126.5
$ws2.Range("H132").Value = 3477.96
$ws2.Range("I132").Value = 3419.9375
$ws2.Range("J132").Value = 3581.111
$ws2.Range("K132").Value = 10259.8125
$ws2.Range("L132").Value = 10743.333
$ws2.Range("M132").Value = -7729.8125
$ws2.Range("N132").Value = -15803.333

# --- BSM ---
$ws3.Range("H4").Value = 262
$ws3.Range("I4").Value = 190
$ws3.Range("J4").Value = 550
$ws3.Range("K4").Value = 190
$ws3.Range("L4").Value = 550
$ws3.Range("M4").Value = -75
$ws3.Range("N4").Value = -780
$ws3.Range("H132").Value = 41500
$ws3.Range("J132").Value = 41500
$ws3.Range("L132").Value = 41500
$ws3.Range("N132").Value = -51620

# --- CRP ---
$ws4.Range("H31").Value = 2319.1
$ws4.Range("I31").Value = 2031.8182
$ws4.Range("J31").Value = 3673.4285
$ws4.Range("K31").Value = 2031.8182
$ws4.Range("L31").Value = 3673.4285
$ws4.Range("M31").Value = -1736.8182
$ws4.Range("N31").Value = -4263.4285
$ws4.Range("H34").Value = 2319.1
$ws4.Range("I34").Value = 2031.8182
$ws4.Range("J34").Value = 3673.4285
$ws4.Range("K34").Value = 2031.8182
$ws4.Range("L34").Value = 3673.4285
$ws4.Range("M34").Value = -1829.8182
$ws4.Range("N34").Value = -4077.4285
$ws4.Range("H122").Value = 1306.6897
$ws4.Range("I122").Value = 1198.381
$ws4.Range("J122").Value = 1591
$ws4.Range("K122").Value = 3595.143
$ws4.Range("L122").Value = 4773
$ws4.Range("M122").Value = -1145.143
$ws4.Range("N122").Value = -9673

# --- CUL ---
$ws5.Range("H5").Value = 420.08334
$ws5.Range("I5").Value = 428.27274
$ws5.Range("J5").Value = 330
$ws5.Range("K5").Value = 1284.81822
$ws5.Range("L5").Value = 990
$ws5.Range("M5").Value = -1172.81822
$ws5.Range("N5").Value = -1214
$ws5.Range("H109").Value = 3391.4119
$ws5.Range("I109").Value = 1325.6666
$ws5.Range("J109").Value = 4518.1816
$ws5.Range("K109").Value = 3976.9998
$ws5.Range("L109").Value = 13554.5448
$ws5.Range("M109").Value = -2936.9998
$ws5.Range("N109").Value = -15634.5448
$ws5.Range("H112").Value = 93866.63
$ws5.Range("I112").Value = 1834.3334
$ws5.Range("J112").Value = 128378.75
$ws5.Range("K112").Value = 5503.0002
$ws5.Range("L112").Value = 385136.25
$ws5.Range("M112").Value = -4395.0002
$ws5.Range("N112").Value = -387352.25
$ws5.Range("H113").Value = 1379739.5
$ws5.Range("I113").Value = 2873900.8
$ws5.Range("J113").Value = 513.7692
$ws5.Range("K113").Value = 8621702.399999999
$ws5.Range("L113").Value = 1541.3076
$ws5.Range("M113").Value = -8619532.399999999
$ws5.Range("N113").Value = -5881.3076
$ws5.Range("H114").Value = 1870
$ws5.Range("I114").Value = 1317.8889
$ws5.Range("J114").Value = 2252.2307
$ws5.Range("K114").Value = 3953.6667
$ws5.Range("L114").Value = 6756.6921
$ws5.Range("M114").Value = -699.6666999999998
$ws5.Range("N114").Value = -13264.6921
$ws5.Range("H115").Value = 2367.6365
$ws5.Range("I115").Value = 1133.8
$ws5.Range("J115").Value = 3395.8333
$ws5.Range("K115").Value = 3401.4
$ws5.Range("L115").Value = 10187.4999
$ws5.Range("M115").Value = -2226.4
$ws5.Range("N115").Value = -12537.4999
$ws5.Range("H116").Value = 1993.5385
$ws5.Range("I116").Value = 1226.8572
$ws5.Range("J116").Value = 2888
$ws5.Range("K116").Value = 3680.5716
$ws5.Range("L116").Value = 8664
$ws5.Range("M116").Value = -238.5715999999998
$ws5.Range("N116").Value = -15548
$ws5.Range("H117").Value = 2958.8333
$ws5.Range("J117").Value = 2940.9285
$ws5.Range("L117").Value = 8822.7855
$ws5.Range("N117").Value = -15706.7855
$ws5.Range("H118").Value = 1560.6666
$ws5.Range("I118").Value = 432.7143
$ws5.Range("J118").Value = 3139.8
$ws5.Range("K118").Value = 1298.1429
$ws5.Range("L118").Value = 9419.400000000001
$ws5.Range("M118").Value = -55.14289999999983
$ws5.Range("N118").Value = -11905.4
$ws5.Range("H119").Value = 3304.2856
$ws5.Range("I119").Value = 2306
$ws5.Range("K119").Value = 6918
$ws5.Range("M119").Value = -2080
$ws5.Range("H120").Value = 8779.125
$ws5.Range("I120").Value = 5875
$ws5.Range("J120").Value = 11683.25
$ws5.Range("K120").Value = 17625
$ws5.Range("L120").Value = 35049.75
$ws5.Range("M120").Value = -12787
$ws5.Range("N120").Value = -44725.75
$ws5.Range("H121").Value = 33334528
$ws5.Range("I121").Value = 550
$ws5.Range("J121").Value = 50001516
$ws5.Range("K121").Value = 1650
$ws5.Range("L121").Value = 150004548
$ws5.Range("M121").Value = -340
$ws5.Range("N121").Value = -150007168
$ws5.Range("H122").Value = 11950.526
$ws5.Range("I122").Value = 21766.1
$ws5.Range("J122").Value = 1044.3334
$ws5.Range("K122").Value = 195894.9
$ws5.Range("L122").Value = 9399.000599999999
$ws5.Range("M122").Value = -193444.9
$ws5.Range("N122").Value = -14299.0006
$ws5.Range("H131").Value = 714.75
$ws5.Range("J131").Value = 945.5714
$ws5.Range("L131").Value = 2836.7142
$ws5.Range("N131").Value = -12916.7142
$ws5.Range("H135").Value = 420.08334
$ws5.Range("I135").Value = 428.27274
$ws5.Range("J135").Value = 330
$ws5.Range("K135").Value = 3854.45466
$ws5.Range("L135").Value = 2970
$ws5.Range("M135").Value = -1319.45466
$ws5.Range("N135").Value = -8040

# --- GSM ---
$ws6.Range("H113").Value = 7526.6313
$ws6.Range("I113").Value = 1849.7778
$ws6.Range("J113").Value = 12635.8
$ws6.Range("K113").Value = 1849.7778
$ws6.Range("L113").Value = 12635.8
$ws6.Range("M113").Value = 320.2221999999999
$ws6.Range("N113").Value = -16975.8

# --- LTW ---
$ws7.Range("H127").Value = 39999.5
$ws7.Range("J127").Value = 39999.5
$ws7.Range("L127").Value = 39999.5
$ws7.Range("N127").Value = -49919.5
$ws7.Range("H136").Value = 0
$ws7.Range("I136").Value = 0
$ws7.Range("J136").Value = 0
$ws7.Range("K136").Value = 0
$ws7.Range("L136").Value = 0
$ws7.Range("M136").ClearContents()
$ws7.Range("N136").ClearContents()

# --- WVR ---
$ws8.Range("H46").Value = 32880
$ws8.Range("J46").Value = 32880
$ws8.Range("L46").Value = 32880
$ws8.Range("N46").Value = -33342
$ws8.Range("H134").Value = 32880
$ws8.Range("J134").Value = 32880
$ws8.Range("L134").Value = 98640
$ws8.Range("N134").Value = -103710
$ws8.Range("H140").Value = 73419.5
$ws8.Range("J140").Value = 73419.5
$ws8.Range("L140").Value = 73419.5
$ws8.Range("N140").Value = -83779.5
